$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 37, pushing the existing rows 37:91 down to 38:92.
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new record.
$ws.Cells.Item(37, 1).Value = 7
$ws.Cells.Item(37, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(37, 3).Value = "Ñuble"
$ws.Cells.Item(37, 4).Value = 44791
$ws.Cells.Item(37, 5).Value = 16
$ws.Cells.Item(37, 6).Value = 100112021
$ws.Cells.Item(37, 7).Value = "Ají"
$ws.Cells.Item(37, 8).Value = "Inferno"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 60
$ws.Cells.Item(37, 11).Value = 17000
$ws.Cells.Item(37, 12).Value = 18000
$ws.Cells.Item(37, 13).Value = 17500
$ws.Cells.Item(37, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(37, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(37, 16).Value = 1167
$ws.Cells.Item(37, 17).Value = 15
$ws.Cells.Item(37, 18).Value = "Hortaliza"
